# Update "想去人数" (interest count) values across sheets, matching the
# output regenerated by the site's data-refresh workflow.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 13023
$ws1.Range("F4").Value  = 30
$ws1.Range("F5").Value  = 85
$ws1.Range("F6").Value  = 95
$ws1.Range("F7").Value  = 54
$ws1.Range("F10").Value = 12994
$ws1.Range("F12").Value = 45
$ws1.Range("F13").Value = 8723
$ws1.Range("F14").Value = 7746
$ws1.Range("F15").Value = 207
$ws1.Range("F19").Value = 991
$ws1.Range("F22").Value = 383
$ws1.Range("F23").Value = 186
$ws1.Range("F24").Value = 333
$ws1.Range("F25").Value = 89

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 4

# --- Sheet "全部类型" (All types, combined listing) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 13023
$ws4.Range("F5").Value  = 30
$ws4.Range("F6").Value  = 85
$ws4.Range("F7").Value  = 95
$ws4.Range("F8").Value  = 54
$ws4.Range("F11").Value = 12994
$ws4.Range("F13").Value = 45
$ws4.Range("F14").Value = 8723
$ws4.Range("F15").Value = 7746
$ws4.Range("F16").Value = 207
$ws4.Range("F20").Value = 991
$ws4.Range("F23").Value = 4
$ws4.Range("F25").Value = 383
$ws4.Range("F26").Value = 186
$ws4.Range("F27").Value = 333
$ws4.Range("F28").Value = 89
